# Generate Report for Handoff
# Update the "Latest Handoff Datetime" (column D) for the ff728b56... file
# row (row 5) on both the zh-cn and de-de localization status sheets to
# reflect a new handoff that was just generated.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-03-08 08:05:01"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-03-08 08:05:14"
